$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Optimistic")
$arr = New-Object 'object[,]' 10,28
$arr[0,0] = -3.062472067885806
$arr[0,1] = -1.238964260942972
$arr[0,2] = 7.437373840288144
$arr[0,3] = -10.72152501685048
$arr[0,4] = -6.291011590403967
$arr[0,5] = -5.122548905943133
$arr[0,6] = -12.09360305353633
$arr[0,7] = -2.194374853749748
$arr[0,8] = 5.588484447958806
$arr[0,9] = 2.942436717004137
$arr[0,10] = 16.97104816942395
$arr[0,11] = 14.66687257116591
$arr[0,12] = -5.279801303487293
$arr[0,13] = -2.027855489921127
$arr[0,14] = -6.020561393380905
$arr[0,15] = -7.19228927023893
$arr[0,16] = -10.88254097223103
$arr[0,17] = -1.91471115892682
$arr[0,18] = -13.10694981157039
$arr[0,19] = -17.32416129087358
$arr[0,20] = -6.89214987764487
$arr[0,21] = -21.53375914133988
$arr[0,22] = -15.59812707021529
$arr[0,23] = -23.40566861543854
$arr[0,24] = -30.63856331673106
$arr[0,25] = -30.35127888996791
$arr[0,26] = -37.94820144614945
$arr[0,27] = -53.50836577097834
$arr[1,0] = -1.139485085637585
$arr[1,1] = -1.241979489516393
$arr[1,2] = 5.59988776736893
$arr[1,3] = -10.20867042735266
$arr[1,4] = -5.327468593538963
$arr[1,5] = -4.685846710856667
$arr[1,6] = -12.10355983859855
$arr[1,7] = 0.4370724523983638
$arr[1,8] = 6.602389756637396
$arr[1,9] = 3.73543412644946
$arr[1,10] = 14.68014145885492
$arr[1,11] = 11.86278066764788
$arr[1,12] = -8.006074568606715
$arr[1,13] = -4.177002472873207
$arr[1,14] = -5.132740790948361
$arr[1,15] = -5.061576974475349
$arr[1,16] = -10.13486565098249
$arr[1,17] = -1.994121656626827
$arr[1,18] = -14.97482099342375
$arr[1,19] = -17.25606074705173
$arr[1,20] = -6.01774991476924
$arr[1,21] = -21.98781703417996
$arr[1,22] = -14.28049012570534
$arr[1,23] = -24.60085862736441
$arr[1,24] = -29.77937129340319
$arr[1,25] = -30.99139082150541
$arr[1,26] = -36.71630962320474
$arr[1,27] = -54.51464059041903
$arr[2,0] = -0.4467285704480823
$arr[2,1] = -1.168513231793739
$arr[2,2] = 7.697884648782734
$arr[2,3] = -11.45940576018544
$arr[2,4] = -7.420421936007324
$arr[2,5] = -4.555631857672106
$arr[2,6] = -11.49305199086962
$arr[2,7] = -0.5261801105283777
$arr[2,8] = 6.537566657249014
$arr[2,9] = 3.856865919024551
$arr[2,10] = 15.27426585997895
$arr[2,11] = 13.30410813477789
$arr[2,12] = -6.82700796225998
$arr[2,13] = -5.75984944657945
$arr[2,14] = -4.370530604408938
$arr[2,15] = -6.163529070313156
$arr[2,16] = -11.39575175159627
$arr[2,17] = -1.23420965526558
$arr[2,18] = -16.14478264679504
$arr[2,19] = -16.51875293779776
$arr[2,20] = -7.536305029284559
$arr[2,21] = -20.97968121432941
$arr[2,22] = -16.3546689063688
$arr[2,23] = -24.26319529170434
$arr[2,24] = -29.41382355194741
$arr[2,25] = -31.29069855197784
$arr[2,26] = -36.49788203024799
$arr[2,27] = -53.03119749409934
$arr[3,0] = -0.5429745515566462
$arr[3,1] = -1.907147413067737
$arr[3,2] = 7.905499311170567
$arr[3,3] = -11.09177154325076
$arr[3,4] = -7.26827331311563
$arr[3,5] = -3.315373223663929
$arr[3,6] = -13.38844145200353
$arr[3,7] = 0.2039540758508949
$arr[3,8] = 7.914190018704072
$arr[3,9] = 3.255623552035985
$arr[3,10] = 16.73886608065255
$arr[3,11] = 12.7294817630945
$arr[3,12] = -7.588763796178019
$arr[3,13] = -7.31327213910574
$arr[3,14] = -4.935316698190776
$arr[3,15] = -6.431944139493824
$arr[3,16] = -11.72007626584621
$arr[3,17] = -1.314255805378605
$arr[3,18] = -15.72796925759258
$arr[3,19] = -18.32290296863097
$arr[3,20] = -5.990521356331409
$arr[3,21] = -20.6133066635805
$arr[3,22] = -14.0550277387915
$arr[3,23] = -23.99766606839822
$arr[3,24] = -29.32499457967619
$arr[3,25] = -30.48579605173391
$arr[3,26] = -37.78510372650722
$arr[3,27] = -53.96948521671438
$arr[4,0] = -1.073598996755198
$arr[4,1] = -2.115623866512848
$arr[4,2] = 3.526017737593814
$arr[4,3] = -11.72149363694911
$arr[4,4] = -4.575912535234156
$arr[4,5] = -4.521514667102897
$arr[4,6] = -12.08755308252102
$arr[4,7] = -1.296319510479982
$arr[4,8] = 7.627759666067703
$arr[4,9] = 4.571844691831569
$arr[4,10] = 15.32272888177511
$arr[4,11] = 13.50073994979975
$arr[4,12] = -6.822117422162766
$arr[4,13] = -4.953462020505739
$arr[4,14] = -3.013847306722354
$arr[4,15] = -4.974567279080247
$arr[4,16] = -11.20396251754573
$arr[4,17] = 1.304574022894279
$arr[4,18] = -16.77690036287866
$arr[4,19] = -17.34588008284184
$arr[4,20] = -7.314957841071942
$arr[4,21] = -21.4264139595598
$arr[4,22] = -16.63070443787229
$arr[4,23] = -22.91160814922423
$arr[4,24] = -28.35527554478023
$arr[4,25] = -30.7576901774996
$arr[4,26] = -37.09499468405851
$arr[4,27] = -55.2068264028877
$arr[5,0] = -1.087523666832824
$arr[5,1] = -2.539149916955957
$arr[5,2] = 6.179486134219745
$arr[5,3] = -11.69205347595221
$arr[5,4] = -6.512786353641327
$arr[5,5] = -6.15082526178726
$arr[5,6] = -12.83287523592911
$arr[5,7] = 0.6290464570105474
$arr[5,8] = 6.836757420075843
$arr[5,9] = 3.310268337863547
$arr[5,10] = 16.46279429234692
$arr[5,11] = 14.28808018590276
$arr[5,12] = -6.174624918148044
$arr[5,13] = -6.193914319411046
$arr[5,14] = -4.432347038087235
$arr[5,15] = -5.339031836287031
$arr[5,16] = -9.712828878304435
$arr[5,17] = 0.5782101881208774
$arr[5,18] = -15.36616464081365
$arr[5,19] = -18.78299958208659
$arr[5,20] = -7.925757249365137
$arr[5,21] = -19.89412215045123
$arr[5,22] = -15.70463466475092
$arr[5,23] = -22.10232708259301
$arr[5,24] = -30.36817266681993
$arr[5,25] = -29.5832863704983
$arr[5,26] = -37.09050088009874
$arr[5,27] = -53.34625897288411
$arr[6,0] = -3.641622573006717
$arr[6,1] = -2.12900126241522
$arr[6,2] = 7.163279691192749
$arr[6,3] = -7.135468866290093
$arr[6,4] = -7.704689008577123
$arr[6,5] = -4.194524896508128
$arr[6,6] = -14.05592443741154
$arr[6,7] = -0.3574066867107408
$arr[6,8] = 5.925438770235912
$arr[6,9] = 4.335126974127458
$arr[6,10] = 12.6836183114694
$arr[6,11] = 10.37611453724694
$arr[6,12] = -6.846003504399961
$arr[6,13] = -5.796429283940403
$arr[6,14] = -4.83798088958136
$arr[6,15] = -7.102074605198542
$arr[6,16] = -11.56399954035802
$arr[6,17] = -0.6784239021190412
$arr[6,18] = -14.69132981903824
$arr[6,19] = -17.35754109112679
$arr[6,20] = -9.480372106978752
$arr[6,21] = -21.90683987664676
$arr[6,22] = -15.67292248391976
$arr[6,23] = -25.55161353390026
$arr[6,24] = -30.45505905808368
$arr[6,25] = -30.50210707585311
$arr[6,26] = -35.75691886645139
$arr[6,27] = -53.27468525189303
$arr[7,0] = -0.469008799770898
$arr[7,1] = -0.1215279288577285
$arr[7,2] = 6.735522609015954
$arr[7,3] = -8.24856719086303
$arr[7,4] = -5.961402250005619
$arr[7,5] = -7.029823983992267
$arr[7,6] = -14.35925005765417
$arr[7,7] = 0.4358911300878008
$arr[7,8] = 7.864173085910287
$arr[7,9] = 3.743470526606636
$arr[7,10] = 15.67926011973778
$arr[7,11] = 12.2806199563457
$arr[7,12] = -6.680354887084606
$arr[7,13] = -4.457464054912472
$arr[7,14] = -3.225687441029581
$arr[7,15] = -6.713225499562106
$arr[7,16] = -12.73490312127528
$arr[7,17] = -0.5900597187203074
$arr[7,18] = -14.20051960605294
$arr[7,19] = -17.45230079500139
$arr[7,20] = -4.501457351254346
$arr[7,21] = -20.84321497597252
$arr[7,22] = -14.92822216379568
$arr[7,23] = -23.46918605518842
$arr[7,24] = -28.76093498514265
$arr[7,25] = -31.26265421115872
$arr[7,26] = -36.39032504734379
$arr[7,27] = -52.92511420596975
$arr[8,0] = -2.080365464378204
$arr[8,1] = -1.90707729942429
$arr[8,2] = 4.023061872667309
$arr[8,3] = -10.41570178435389
$arr[8,4] = -8.727219954383457
$arr[8,5] = -2.83719562236254
$arr[8,6] = -12.38814932361473
$arr[8,7] = -0.6060173421286192
$arr[8,8] = 5.767738540023842
$arr[8,9] = 1.672794019622382
$arr[8,10] = 17.91598939122799
$arr[8,11] = 12.67724040355363
$arr[8,12] = -7.292803488780454
$arr[8,13] = -6.171759807130333
$arr[8,14] = -5.84851643341989
$arr[8,15] = -6.170196315530535
$arr[8,16] = -10.59172153935027
$arr[8,17] = -1.286585587975646
$arr[8,18] = -17.46323388142073
$arr[8,19] = -18.25495878730424
$arr[8,20] = -6.965894922444066
$arr[8,21] = -21.84572175772603
$arr[8,22] = -18.36768470713927
$arr[8,23] = -23.03688081070551
$arr[8,24] = -28.42885024273507
$arr[8,25] = -28.07001013433978
$arr[8,26] = -38.01505237345851
$arr[8,27] = -53.83160187051888
$arr[9,0] = -0.1361167218341284
$arr[9,1] = -2.349945125598151
$arr[9,2] = 5.513238801517774
$arr[9,3] = -11.72874705317787
$arr[9,4] = -6.274310057289685
$arr[9,5] = -4.632148106995988
$arr[9,6] = -11.81772026105496
$arr[9,7] = -0.839997031347131
$arr[9,8] = 7.47240550580891
$arr[9,9] = 2.505299535104602
$arr[9,10] = 16.71846071493313
$arr[9,11] = 13.17896855237374
$arr[9,12] = -6.720551809843916
$arr[9,13] = -6.563113782472557
$arr[9,14] = -6.363248297528925
$arr[9,15] = -6.427359880421676
$arr[9,16] = -11.513231326504
$arr[9,17] = -1.921645304521856
$arr[9,18] = -14.77107448433218
$arr[9,19] = -17.23708430482464
$arr[9,20] = -6.361845134784749
$arr[9,21] = -21.49187943345967
$arr[9,22] = -16.01014149161149
$arr[9,23] = -22.63695321803819
$arr[9,24] = -30.0093366059017
$arr[9,25] = -29.07122773655565
$arr[9,26] = -38.53271419968799
$arr[9,27] = -52.6451091552
$ws.Range("B2:AC11").Value = $arr

$ws = $wb.Worksheets.Item("Pessimistic")
$arr = New-Object 'object[,]' 10,28
$arr[0,0] = 18.09035783948996
$arr[0,1] = -7.574573752029489
$arr[0,2] = -3.971308705849237
$arr[0,3] = 3.749968886309238
$arr[0,4] = -1.910196971052603
$arr[0,5] = -13.51254108104425
$arr[0,6] = 3.506309586628578
$arr[0,7] = -20.79024072154307
$arr[0,8] = -17.57367864959065
$arr[0,9] = -3.189376465185704
$arr[0,10] = 0.9753755475010712
$arr[0,11] = -5.991652528988276
$arr[0,12] = -2.411925080305639
$arr[0,13] = -1.380797086968166
$arr[0,14] = 3.243807476820508
$arr[0,15] = -4.91913646837151
$arr[0,16] = -10.80249251230266
$arr[0,17] = 2.540598671141096
$arr[0,18] = -9.568009608949712
$arr[0,19] = 0.6872512751907638
$arr[0,20] = -5.287375586005513
$arr[0,21] = -9.310085586312773
$arr[0,22] = -5.494719952613302
$arr[0,23] = -3.197462258193638
$arr[0,24] = 13.61456451492145
$arr[0,25] = -0.7684429048442256
$arr[0,26] = 7.625931949875342
$arr[0,27] = -13.96191155022901
$arr[1,0] = 18.4275624159164
$arr[1,1] = -7.980532244518947
$arr[1,2] = -5.048261751771806
$arr[1,3] = 4.438889238237997
$arr[1,4] = -3.512504317599351
$arr[1,5] = -16.70872758107307
$arr[1,6] = 1.558129307951976
$arr[1,7] = -19.28215179376657
$arr[1,8] = -20.55900538763284
$arr[1,9] = -2.744776773067609
$arr[1,10] = 0.8033988523865818
$arr[1,11] = -5.725335759733898
$arr[1,12] = -4.540925844589511
$arr[1,13] = 0.3481166177452195
$arr[1,14] = 1.221004696279663
$arr[1,15] = -5.62119624506836
$arr[1,16] = -10.21089775750867
$arr[1,17] = 1.400511705517468
$arr[1,18] = -8.007337948027942
$arr[1,19] = 1.531361832826273
$arr[1,20] = -7.345168339539679
$arr[1,21] = -11.63954658748583
$arr[1,22] = -3.899057106667333
$arr[1,23] = -2.201712314616215
$arr[1,24] = 15.26195019931927
$arr[1,25] = -0.5164229851871598
$arr[1,26] = 4.602895639801431
$arr[1,27] = -15.58925355648864
$arr[2,0] = 18.23352649795185
$arr[2,1] = -8.035318142607101
$arr[2,2] = -2.858965818053119
$arr[2,3] = 3.470267712009265
$arr[2,4] = -2.776720668844645
$arr[2,5] = -14.38145822788996
$arr[2,6] = 3.408941575653318
$arr[2,7] = -20.42744968004643
$arr[2,8] = -20.35952349386834
$arr[2,9] = -5.107111300313607
$arr[2,10] = -0.7447743638399709
$arr[2,11] = -8.483059343605724
$arr[2,12] = -4.610434900243037
$arr[2,13] = -0.5397915753674927
$arr[2,14] = 1.419432126036976
$arr[2,15] = -5.205415541159568
$arr[2,16] = -9.791737111820172
$arr[2,17] = 1.941701660103416
$arr[2,18] = -8.47560736245363
$arr[2,19] = 3.520830109519716
$arr[2,20] = -7.106882944097684
$arr[2,21] = -10.15768330591325
$arr[2,22] = -4.973940996364544
$arr[2,23] = -3.496039517431503
$arr[2,24] = 14.69142479176936
$arr[2,25] = -1.388361879102415
$arr[2,26] = 4.314078173670611
$arr[2,27] = -15.6641561893458
$arr[3,0] = 16.74824878055296
$arr[3,1] = -8.492725661068999
$arr[3,2] = -1.635892768305034
$arr[3,3] = 2.960242945020367
$arr[3,4] = -2.571103187473979
$arr[3,5] = -13.92781360745312
$arr[3,6] = 0.1654592658020624
$arr[3,7] = -18.85186924589374
$arr[3,8] = -17.73545839215408
$arr[3,9] = -4.541625847341345
$arr[3,10] = 0.3093604562131307
$arr[3,11] = -6.071370914845081
$arr[3,12] = -4.815600973864653
$arr[3,13] = 0.648137119660757
$arr[3,14] = 2.075472845400615
$arr[3,15] = -4.634159287251769
$arr[3,16] = -9.593433320518313
$arr[3,17] = 3.02451871712874
$arr[3,18] = -10.43168587258605
$arr[3,19] = -1.231822178536718
$arr[3,20] = -6.243332422082724
$arr[3,21] = -11.11698676702065
$arr[3,22] = -7.024427682695825
$arr[3,23] = -3.753213486751616
$arr[3,24] = 13.73796377622137
$arr[3,25] = -0.8573603035361135
$arr[3,26] = 4.17797802298842
$arr[3,27] = -12.89907689256381
$arr[4,0] = 18.4401846518541
$arr[4,1] = -6.150587008608293
$arr[4,2] = -3.924958097387672
$arr[4,3] = 3.766617717812798
$arr[4,4] = -3.234672552232434
$arr[4,5] = -15.07776567605433
$arr[4,6] = 3.456250992684136
$arr[4,7] = -20.23267781255048
$arr[4,8] = -19.52335523172497
$arr[4,9] = -3.547447395729943
$arr[4,10] = 1.765053511039139
$arr[4,11] = -6.156080600940907
$arr[4,12] = -2.934395947056335
$arr[4,13] = -0.337270924805789
$arr[4,14] = 1.491518503012746
$arr[4,15] = -7.108608720298049
$arr[4,16] = -9.156337668644035
$arr[4,17] = 3.838896899643307
$arr[4,18] = -7.581885959570799
$arr[4,19] = 0.8524986813955251
$arr[4,20] = -7.563215938437746
$arr[4,21] = -11.9098593723729
$arr[4,22] = -5.465884469499989
$arr[4,23] = -3.093702025499964
$arr[4,24] = 12.73787021212637
$arr[4,25] = 0.3827529485764076
$arr[4,26] = 4.947488522940082
$arr[4,27] = -13.04264888814803
$arr[5,0] = 18.06117833311864
$arr[5,1] = -8.408198491391103
$arr[5,2] = -6.159417579528711
$arr[5,3] = 2.923393891631036
$arr[5,4] = -2.783872803658633
$arr[5,5] = -14.74545104782752
$arr[5,6] = 1.545447635041748
$arr[5,7] = -21.01020852996496
$arr[5,8] = -19.19117798882987
$arr[5,9] = -3.501276163988181
$arr[5,10] = 3.155630785915084
$arr[5,11] = -6.902757443293874
$arr[5,12] = -2.574972608017441
$arr[5,13] = -0.2656109178575035
$arr[5,14] = 4.431589926570801
$arr[5,15] = -4.251344794941256
$arr[5,16] = -9.224125455733265
$arr[5,17] = 2.470427075940603
$arr[5,18] = -9.895353367083
$arr[5,19] = 1.465018456585525
$arr[5,20] = -7.261393535840777
$arr[5,21] = -9.87301489316894
$arr[5,22] = -5.554877043626924
$arr[5,23] = -3.700949458856035
$arr[5,24] = 13.52987950596424
$arr[5,25] = -1.071404630901701
$arr[5,26] = 4.59532044427843
$arr[5,27] = -13.98868896652222
$arr[6,0] = 14.46545461348884
$arr[6,1] = -8.79909200947262
$arr[6,2] = -4.180164659267599
$arr[6,3] = 4.522378902259005
$arr[6,4] = -3.557531608685616
$arr[6,5] = -15.35717753823863
$arr[6,6] = 1.292585346167598
$arr[6,7] = -23.19667014036774
$arr[6,8] = -19.76720161514999
$arr[6,9] = -3.380960565534672
$arr[6,10] = 0.7429976963933389
$arr[6,11] = -5.701349594812908
$arr[6,12] = -4.016779093079797
$arr[6,13] = 0.8096796638140784
$arr[6,14] = 2.250026583761316
$arr[6,15] = -5.87143080077413
$arr[6,16] = -9.793690355868247
$arr[6,17] = 0.7757422195845414
$arr[6,18] = -8.763227521990894
$arr[6,19] = 1.013876613464367
$arr[6,20] = -6.580179954455257
$arr[6,21] = -10.8433298769826
$arr[6,22] = -5.149002375238927
$arr[6,23] = -2.88054145963446
$arr[6,24] = 11.88410578825684
$arr[6,25] = 0.6011017056205619
$arr[6,26] = 3.422149956519359
$arr[6,27] = -15.77374895441607
$arr[7,0] = 17.34309386404055
$arr[7,1] = -10.1879629740546
$arr[7,2] = -3.048813949778206
$arr[7,3] = 4.351995728710372
$arr[7,4] = -1.794783157429226
$arr[7,5] = -13.97176006582244
$arr[7,6] = 2.949519352447306
$arr[7,7] = -19.32613876171793
$arr[7,8] = -17.82339642374711
$arr[7,9] = -3.261762399869323
$arr[7,10] = -0.24647128397533
$arr[7,11] = -6.555383926371476
$arr[7,12] = -5.250840221259693
$arr[7,13] = 0.9194265532273465
$arr[7,14] = 1.955633960221289
$arr[7,15] = -4.124096880647417
$arr[7,16] = -9.585289446091629
$arr[7,17] = 2.594304933755851
$arr[7,18] = -10.75024937324942
$arr[7,19] = 3.4572053456851
$arr[7,20] = -5.412579535193903
$arr[7,21] = -10.46583974295816
$arr[7,22] = -3.937830426850174
$arr[7,23] = -1.473525862779714
$arr[7,24] = 15.17519763217135
$arr[7,25] = 0.6701579828637394
$arr[7,26] = 5.195131980509358
$arr[7,27] = -12.73183277306354
$arr[8,0] = 18.88246445693753
$arr[8,1] = -10.53983468770544
$arr[8,2] = -6.218369579543585
$arr[8,3] = 4.320303276887742
$arr[8,4] = -2.671644800466535
$arr[8,5] = -13.61845552953149
$arr[8,6] = 2.337340272206443
$arr[8,7] = -20.45544284217544
$arr[8,8] = -17.91074860348461
$arr[8,9] = -5.435375523242164
$arr[8,10] = 0.2428245534378091
$arr[8,11] = -7.334362011457928
$arr[8,12] = -4.194095086798195
$arr[8,13] = 0.1262741435411319
$arr[8,14] = 1.924387723579986
$arr[8,15] = -3.245067060336284
$arr[8,16] = -8.415493670390351
$arr[8,17] = 1.702050414269175
$arr[8,18] = -9.373913261778794
$arr[8,19] = 1.569027089378332
$arr[8,20] = -6.070520006175636
$arr[8,21] = -10.03862928418067
$arr[8,22] = -4.37727945013542
$arr[8,23] = -3.373936862306573
$arr[8,24] = 14.20784948955467
$arr[8,25] = -0.973573076573599
$arr[8,26] = 5.044517136517763
$arr[8,27] = -14.18833798152048
$arr[9,0] = 14.39526934208785
$arr[9,1] = -6.970670747275994
$arr[9,2] = -4.023757878161371
$arr[9,3] = 2.745578194932317
$arr[9,4] = -3.277047007142355
$arr[9,5] = -14.01941111081195
$arr[9,6] = 2.191424304184737
$arr[9,7] = -20.36257599291497
$arr[9,8] = -20.62332079988082
$arr[9,9] = -2.946206705165964
$arr[9,10] = -0.1914215563262407
$arr[9,11] = -6.677325763711103
$arr[9,12] = -4.728844610780051
$arr[9,13] = -0.3556223103348304
$arr[9,14] = 1.477631447900896
$arr[9,15] = -5.711796545685422
$arr[9,16] = -10.09359591570766
$arr[9,17] = 1.595580855035731
$arr[9,18] = -9.223945395925874
$arr[9,19] = 1.934153584457774
$arr[9,20] = -6.452304236101971
$arr[9,21] = -12.14759258267853
$arr[9,22] = -4.223342938054854
$arr[9,23] = -2.350766848609481
$arr[9,24] = 12.90649797236187
$arr[9,25] = -0.8837959196796081
$arr[9,26] = 4.72300243636633
$arr[9,27] = -12.72557155931403
$ws.Range("B2:AC11").Value = $arr

$ws = $wb.Worksheets.Item("Middle")
$arr = New-Object 'object[,]' 10,28
$arr[0,0] = -5.625236089776784
$arr[0,1] = -4.492804866806801
$arr[0,2] = -4.713966740853254
$arr[0,3] = 1.203317352989799
$arr[0,4] = -14.06284944644749
$arr[0,5] = -12.78727537594284
$arr[0,6] = 2.814965286016304
$arr[0,7] = -2.53913221659873
$arr[0,8] = -8.563753893333859
$arr[0,9] = 1.488950011511102
$arr[0,10] = -11.38646500702958
$arr[0,11] = 4.423955762620237
$arr[0,12] = -8.128517429177702
$arr[0,13] = -5.240345259732992
$arr[0,14] = -13.03817715404549
$arr[0,15] = -13.78024308438627
$arr[0,16] = -0.7503579798194071
$arr[0,17] = -11.78761853554003
$arr[0,18] = -12.22987088733476
$arr[0,19] = 6.38124098023076
$arr[0,20] = -9.971296061096416
$arr[0,21] = -7.901048654200913
$arr[0,22] = -16.3878612140004
$arr[0,23] = -1.828921386879244
$arr[0,24] = -6.595812862720742
$arr[0,25] = -27.39268680411794
$arr[0,26] = -6.291168289940297
$arr[0,27] = -19.80186137763567
$arr[1,0] = -6.818883773261904
$arr[1,1] = -5.372510623441075
$arr[1,2] = -2.368929407152529
$arr[1,3] = 0.5252353418916407
$arr[1,4] = -15.00014171096014
$arr[1,5] = -11.20677851816033
$arr[1,6] = 3.272069460222101
$arr[1,7] = -3.390386788649379
$arr[1,8] = -9.661424112147625
$arr[1,9] = 3.563205879808467
$arr[1,10] = -10.48280988227735
$arr[1,11] = 4.051915199238149
$arr[1,12] = -9.114230651989303
$arr[1,13] = -6.436374722914955
$arr[1,14] = -13.2941314263707
$arr[1,15] = -13.10836807070307
$arr[1,16] = -1.053363912535194
$arr[1,17] = -12.18358353125703
$arr[1,18] = -11.75816231033135
$arr[1,19] = 3.948180534477614
$arr[1,20] = -11.24216040821767
$arr[1,21] = -6.205397729379446
$arr[1,22] = -15.35892040533156
$arr[1,23] = -3.221889975387329
$arr[1,24] = -7.58937915628077
$arr[1,25] = -26.61643965109187
$arr[1,26] = -9.798801387419182
$arr[1,27] = -18.33582224848669
$arr[2,0] = -6.577149976098827
$arr[2,1] = -4.84254589086262
$arr[2,2] = -4.126286625305527
$arr[2,3] = 1.426689940638771
$arr[2,4] = -13.01651984368936
$arr[2,5] = -9.92046812818895
$arr[2,6] = 2.707003519207431
$arr[2,7] = -2.281426627988874
$arr[2,8] = -10.68827321018486
$arr[2,9] = 3.109659809836429
$arr[2,10] = -7.814209627484809
$arr[2,11] = 3.339333410705613
$arr[2,12] = -7.861108158834494
$arr[2,13] = -4.824517096675294
$arr[2,14] = -11.46419735514247
$arr[2,15] = -13.30110365178605
$arr[2,16] = 0.1951328909401382
$arr[2,17] = -12.26144680907025
$arr[2,18] = -11.96877888669006
$arr[2,19] = 7.239012983944118
$arr[2,20] = -10.99581300095454
$arr[2,21] = -6.517114725015108
$arr[2,22] = -17.53393828335611
$arr[2,23] = -4.565812851299319
$arr[2,24] = -6.187734418891328
$arr[2,25] = -25.40702851653849
$arr[2,26] = -8.221784165018422
$arr[2,27] = -19.53115493713051
$arr[3,0] = -7.936111302687825
$arr[3,1] = -4.895449169464842
$arr[3,2] = -3.263112850524579
$arr[3,3] = 1.355114777011556
$arr[3,4] = -12.81354458248061
$arr[3,5] = -11.67688882982877
$arr[3,6] = 2.090104280098602
$arr[3,7] = -2.921948632120482
$arr[3,8] = -9.26592318707321
$arr[3,9] = 2.396944602740177
$arr[3,10] = -7.886432883505679
$arr[3,11] = 2.281685840803264
$arr[3,12] = -8.441932340827353
$arr[3,13] = -6.186082268700011
$arr[3,14] = -12.23152380444628
$arr[3,15] = -12.50049947678504
$arr[3,16] = -0.6535831533651897
$arr[3,17] = -12.30761874478687
$arr[3,18] = -12.1831014301878
$arr[3,19] = 3.975848884848237
$arr[3,20] = -10.61156376686587
$arr[3,21] = -7.658041451829229
$arr[3,22] = -16.44785058316933
$arr[3,23] = -1.221050464013117
$arr[3,24] = -8.836200761572334
$arr[3,25] = -27.68346271763291
$arr[3,26] = -8.664621425219382
$arr[3,27] = -19.08343402808148
$arr[4,0] = -7.08567176177534
$arr[4,1] = -2.649540344830688
$arr[4,2] = -2.255586109612517
$arr[4,3] = 0.725118117428257
$arr[4,4] = -13.21249038820341
$arr[4,5] = -11.66417873205729
$arr[4,6] = 0.8273670033958274
$arr[4,7] = -2.786073701134428
$arr[4,8] = -9.543148242840068
$arr[4,9] = 2.529770891546195
$arr[4,10] = -10.3836400872925
$arr[4,11] = 4.875378230872986
$arr[4,12] = -9.430572262807392
$arr[4,13] = -5.455606548103056
$arr[4,14] = -10.55460815651818
$arr[4,15] = -13.44564667120729
$arr[4,16] = 0.1601442658475829
$arr[4,17] = -11.85373722064415
$arr[4,18] = -13.56888587427268
$arr[4,19] = 3.317018701188685
$arr[4,20] = -10.21551454103358
$arr[4,21] = -6.379514119553577
$arr[4,22] = -16.91442713884609
$arr[4,23] = -3.436446022981291
$arr[4,24] = -6.828959019404071
$arr[4,25] = -26.0516026812103
$arr[4,26] = -7.382764177251795
$arr[4,27] = -20.49449588931284
$arr[5,0] = -8.615125846761462
$arr[5,1] = -6.155545757761051
$arr[5,2] = -4.057473078151852
$arr[5,3] = 2.829511413657539
$arr[5,4] = -15.36436145665334
$arr[5,5] = -12.14151816072064
$arr[5,6] = 1.737221146982148
$arr[5,7] = -3.56499233540253
$arr[5,8] = -9.112759131834185
$arr[5,9] = 3.057691617038121
$arr[5,10] = -12.58360150463224
$arr[5,11] = 2.700168831503466
$arr[5,12] = -7.914273274860168
$arr[5,13] = -5.062212063698671
$arr[5,14] = -12.21032190345417
$arr[5,15] = -13.45606897253007
$arr[5,16] = -1.161635980390625
$arr[5,17] = -12.70362659622453
$arr[5,18] = -14.9953192344393
$arr[5,19] = 6.923646561016074
$arr[5,20] = -12.6267630920654
$arr[5,21] = -5.771769279628653
$arr[5,22] = -16.63329244777516
$arr[5,23] = -3.210379950150933
$arr[5,24] = -5.905990818943131
$arr[5,25] = -27.47546556318177
$arr[5,26] = -9.046575940021771
$arr[5,27] = -20.21670086955524
$arr[6,0] = -8.299926586843922
$arr[6,1] = -6.827812898721437
$arr[6,2] = -2.032454541942204
$arr[6,3] = -0.2974517181959255
$arr[6,4] = -12.84065354440332
$arr[6,5] = -13.99021736801504
$arr[6,6] = 2.158945582377878
$arr[6,7] = -2.013872113900789
$arr[6,8] = -9.937641874022187
$arr[6,9] = 1.582005331774441
$arr[6,10] = -9.329658394848437
$arr[6,11] = 4.333857759329925
$arr[6,12] = -8.523535336970319
$arr[6,13] = -5.934677989867392
$arr[6,14] = -11.66751092705597
$arr[6,15] = -12.87920393085149
$arr[6,16] = -0.4664457940986719
$arr[6,17] = -10.80277033875112
$arr[6,18] = -13.25958454342712
$arr[6,19] = 6.225404415846988
$arr[6,20] = -9.958539160250584
$arr[6,21] = -6.249972697109444
$arr[6,22] = -18.52617701423971
$arr[6,23] = -2.545551244676546
$arr[6,24] = -6.318911461166802
$arr[6,25] = -27.18173510509245
$arr[6,26] = -8.717540819089944
$arr[6,27] = -19.43647621383438
$arr[7,0] = -5.159724254257917
$arr[7,1] = -3.608080989648451
$arr[7,2] = -4.546301305201446
$arr[7,3] = 0.5956319669226717
$arr[7,4] = -14.40655854851215
$arr[7,5] = -11.75613910617066
$arr[7,6] = 2.293709898305701
$arr[7,7] = -2.019216740271562
$arr[7,8] = -9.293637565311618
$arr[7,9] = 2.565029370036655
$arr[7,10] = -9.210916492149646
$arr[7,11] = 4.872506153600864
$arr[7,12] = -8.345321430992163
$arr[7,13] = -4.341359969311836
$arr[7,14] = -11.87797614267347
$arr[7,15] = -13.51525027468311
$arr[7,16] = -1.14232425020868
$arr[7,17] = -11.79853450039488
$arr[7,18] = -15.0145481397347
$arr[7,19] = 4.987356892536513
$arr[7,20] = -9.197067949227812
$arr[7,21] = -7.229536028097698
$arr[7,22] = -17.59005194518193
$arr[7,23] = -2.564118925033186
$arr[7,24] = -6.527505775908338
$arr[7,25] = -25.94319784955529
$arr[7,26] = -7.048301767423446
$arr[7,27] = -17.47778807908902
$arr[8,0] = -6.471966992872805
$arr[8,1] = -4.777533999153452
$arr[8,2] = -4.582481203173075
$arr[8,3] = 0.6312624133700635
$arr[8,4] = -13.54413194673393
$arr[8,5] = -11.89163223649107
$arr[8,6] = 2.865124542786676
$arr[8,7] = -0.8684818442268925
$arr[8,8] = -10.23667331526461
$arr[8,9] = 0.7410905638845309
$arr[8,10] = -9.589653547399807
$arr[8,11] = 3.532865117008764
$arr[8,12] = -8.306091898875096
$arr[8,13] = -5.041174295482772
$arr[8,14] = -11.56289191152655
$arr[8,15] = -14.20303666383636
$arr[8,16] = 0.1648068711340822
$arr[8,17] = -12.25564979648634
$arr[8,18] = -11.72788361132882
$arr[8,19] = 7.453913935019592
$arr[8,20] = -11.72819620821181
$arr[8,21] = -6.940399847550059
$arr[8,22] = -16.48626042303003
$arr[8,23] = -2.950341180881527
$arr[8,24] = -5.496251044652269
$arr[8,25] = -25.72884543653182
$arr[8,26] = -7.597157904816182
$arr[8,27] = -18.20968207724654
$arr[9,0] = -7.32011811125771
$arr[9,1] = -4.368631950442298
$arr[9,2] = -1.013048958529922
$arr[9,3] = 2.256965039268387
$arr[9,4] = -13.7437638551725
$arr[9,5] = -11.90884766303352
$arr[9,6] = 4.280563341720634
$arr[9,7] = -2.481826171198228
$arr[9,8] = -9.783734844454637
$arr[9,9] = 1.589096433140175
$arr[9,10] = -10.18433075460696
$arr[9,11] = 3.419682432399803
$arr[9,12] = -8.948829860060492
$arr[9,13] = -5.410170774906486
$arr[9,14] = -12.94768789464042
$arr[9,15] = -13.06812452867731
$arr[9,16] = 0.9350412743278582
$arr[9,17] = -12.96605832212864
$arr[9,18] = -12.79596643688225
$arr[9,19] = 5.597263014028109
$arr[9,20] = -10.87707380483311
$arr[9,21] = -7.501051635170676
$arr[9,22] = -18.15398984148175
$arr[9,23] = -2.219918516582767
$arr[9,24] = -5.163208719552449
$arr[9,25] = -25.80346850010205
$arr[9,26] = -7.785552887824117
$arr[9,27] = -17.48884074992373
$ws.Range("B2:AC11").Value = $arr
